$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

function Set-PlainCell($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

Set-PlainCell "D2" "26.418.45"
Set-PlainCell "E2" "  +1.11%  "
Set-PlainCell "D3" "1.672.78"
Set-PlainCell "E3" "  +1.10%  "
Set-TextCell "D4" "1.011"
Set-PlainCell "E4" "  +0.84%  "
Set-TextCell "D5" "220.96"
Set-PlainCell "E5" "  +1.48%  "
Set-TextCell "D6" "0.5364"
Set-PlainCell "E6" "  +1.26%  "
Set-PlainCell "E7" "  +0.75%  "
Set-TextCell "D8" "0.2671"
Set-PlainCell "E8" "  +2.25%  "
Set-TextCell "D9" "0.06414"
Set-PlainCell "E9" "  +1.38%  "
Set-PlainCell "E10" "  +3.33%  "
Set-TextCell "D11" "0.07866"
Set-PlainCell "E11" "  +0.85%  "
Set-TextCell "D12" "4.578"
Set-PlainCell "E12" "  +1.21%  "
Set-PlainCell "D13" "1.675.56"
Set-PlainCell "E13" "  +2.23%  "
Set-PlainCell "D14" "1.902.02"
Set-PlainCell "E14" "  +1.09%  "
Set-TextCell "D15" "0.5653"
Set-PlainCell "E15" "  +3.04%  "
Set-PlainCell "D16" "0.0₅8209"
Set-PlainCell "E16" "  -0.03%  "
Set-PlainCell "E17" "  +1.56%  "
Set-PlainCell "D18" "26.459.85"
Set-PlainCell "E18" "  +1.31%  "
Set-TextCell "D20" "4.717"
Set-PlainCell "E20" "  +2.68%  "
Set-TextCell "D21" "196.87"
Set-PlainCell "E21" "  +3.04%  "
Set-TextCell "D22" "10.36"
Set-PlainCell "E22" "  +2.84%  "
Set-TextCell "D23" "6.071"
Set-PlainCell "E23" "  +0.79%  "
Set-TextCell "D24" "1.012"
Set-PlainCell "E24" "  +0.78%  "
Set-TextCell "D25" "146.25"
Set-PlainCell "E25" "  +0.86%  "
Set-PlainCell "E26" "  +0.55%  "
Set-TextCell "D27" "7.272"
Set-PlainCell "E27" "  +0.74%  "
Set-TextCell "D28" "16.29"
Set-PlainCell "E28" "  +1.87%  "
Set-TextCell "D29" "1.512"
Set-PlainCell "E29" "  +3.64%  "
Set-TextCell "D30" "0.05909"
Set-PlainCell "E30" "  +2.03%  "
Set-PlainCell "E31" "  +1.60%  "
Set-TextCell "D32" "3.591"
Set-PlainCell "E32" "  +1.19%  "
Set-TextCell "D33" "3.313"
Set-PlainCell "E33" "  +1.34%  "
Set-TextCell "D34" "1.629"
Set-PlainCell "E34" "  +1.94%  "
Set-TextCell "D35" "0.9738"
Set-PlainCell "E35" "  +2.70%  "
Set-TextCell "D36" "2.852"
Set-PlainCell "E36" "  +1.79%  "
Set-PlainCell "E37" "  +0.64%  "
Set-TextCell "D38" "0.5845"
Set-PlainCell "E38" "  +1.49%  "
Set-TextCell "D39" "0.01614"
Set-PlainCell "E39" "  +0.15%  "
Set-PlainCell "D40" "1.078.05"
Set-PlainCell "E40" "  +4.29%  "
Set-PlainCell "B41" "TrustWalletToken"
Set-PlainCell "C41" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell "D41" "0.8701"
Set-PlainCell "E41" "  +1.70%  "
Set-PlainCell "B42" "FraxShare"
Set-PlainCell "C42" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D42" "5.910"
Set-PlainCell "E42" "  +2.83%  "
Set-PlainCell "E43" "  +0.79%  "
Set-TextCell "D44" "104.63"
Set-PlainCell "E44" "  +0.16%  "
Set-PlainCell "D45" "1.811.57"
Set-TextCell "D46" "58.41"
Set-PlainCell "E46" "  +2.55%  "
Set-PlainCell "D47" "0.0₈106"
Set-PlainCell "E47" "  -4.72%  "
Set-TextCell "D48" "1.013"
Set-PlainCell "E48" "  +0.63%  "
Set-TextCell "D49" "0.4402"
Set-TextCell "D50" "8.084"
Set-PlainCell "E50" "  +3.03%  "
Set-TextCell "D51" "0.05169"
Set-PlainCell "E51" "  +0.47%  "
